$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.571.57'
$ws.Range("E2").Value = '  +2.77%  '
$ws.Range("D3").Value = '1.669.73'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4776'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  +2.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06171'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.58%  '
$ws.Range("D10").Value = '1.673.94'
$ws.Range("E10").Value = '  +2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06986'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5881'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.374'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9998'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '25.568.78'
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006746'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.56%  '
$ws.Range("E20").Value = '  +3.64%  '
$ws.Range("D21").Value = '1.886.44'
$ws.Range("E21").Value = '  +2.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.443'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.788'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.258'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.386'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.719'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.94%  '
$ws.Range("E30").Value = '  +6.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07869'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.629'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9989'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04298'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.619'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9542'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6051'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.583'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9218'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9996'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.857'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01474'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3759'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.878'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1118'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("E47").Value = '  +3.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05264'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.433'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.22%  '
